# Weekly update: a new price observation was added for
# "Vega Modelo de Temuco - Achicoria" above the existing rows.
# This inserts one new data row at row 125 (shifting the existing
# rows 125:151 down to 126:152) and populates it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 125, pushing existing data down.
$ws.Rows.Item(125).Insert()

# Fill in the new row with the latest observation.
$ws.Range("A125").Value = 10
$ws.Range("B125").Value = "Vega Modelo de Temuco"
$ws.Range("C125").Value = "La Araucanía"
$ws.Range("D125").Value = 45204
$ws.Range("E125").Value = 9
$ws.Range("F125").Value = 100112010
$ws.Range("G125").Value = "Achicoria"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 65
$ws.Range("K125").Value = 12000
$ws.Range("L125").Value = 12000
$ws.Range("M125").Value = 12000
$ws.Range("N125").Value = "$/caja 18 unidades"
$ws.Range("O125").Value = "Región del Maule"
$ws.Range("P125").Value = 667
$ws.Range("Q125").Value = 18
$ws.Range("R125").Value = "Hortaliza"
